# Regenerate the handback status report: two source files were re-run
# through the pipeline, producing new GUIDs/hashes/timestamps.
#
# file 1: 196f4342-cf7e-4c67-b105-f0f976a4b81c -> c075ce48-2946-4b29-9cd8-dc9e82ee0cbd
# file 2: 650b5967-7ade-4182-9ac8-804e3ebc3ae8 -> ffff8f690006-79ab-4302-987b-ffa063e2dc7b
# xlf content hash (both languages, both files now share it):
#         cc8d7be5a865dc4e067b39464ada1f9b9e2f8dc7 / 918378bc386f3928fb7cab8fb055ddcf83834c05
#           -> e4e5a9398226a09db97282e15f3993711e8ca2fd
# handoff/handback timestamps move forward ~1 minute.

$wb = $excel.ActiveWorkbook

$oldUuid1 = "196f4342-cf7e-4c67-b105-f0f976a4b81c"
$newUuid1 = "c075ce48-2946-4b29-9cd8-dc9e82ee0cbd"
$oldUuid2 = "650b5967-7ade-4182-9ac8-804e3ebc3ae8"
$newUuid2 = "ffff8f690006-79ab-4302-987b-ffa063e2dc7b"

$oldHash1 = "cc8d7be5a865dc4e067b39464ada1f9b9e2f8dc7"
$oldHash2 = "918378bc386f3928fb7cab8fb055ddcf83834c05"
$newHash  = "e4e5a9398226a09db97282e15f3993711e8ca2fd"

# Text substitutions to apply verbatim to every cell value and hyperlink
# display string in the workbook, in order.
$replacements = @{}
$replacements["$oldUuid1.md"] = "$newUuid1.md"
$replacements["$oldUuid2.md"] = "$newUuid2.md"
$replacements["$oldUuid1.$oldHash1.zh-cn.xlf"] = "$newUuid1.$newHash.zh-cn.xlf"
$replacements["$oldUuid2.$oldHash2.zh-cn.xlf"] = "$newUuid1.$newHash.zh-cn.xlf"
$replacements["$oldUuid1.$oldHash1.de-de.xlf"] = "$newUuid1.$newHash.de-de.xlf"
$replacements["$oldUuid2.$oldHash2.de-de.xlf"] = "$newUuid1.$newHash.de-de.xlf"
$replacements["2016-03-20 22:54:16"] = "2016-03-20 22:55:40"
$replacements["2016-03-20 22:54:37"] = "2016-03-20 22:56:01"
$replacements["2016-03-20 22:54:19"] = "2016-03-20 22:55:44"
$replacements["2016-03-20 22:54:42"] = "2016-03-20 22:56:07"

foreach ($ws in $wb.Worksheets) {
    # NOTE: iterating `$used.Rows` / `.Cells` with `foreach` only yields the
    # first column here, so walk the used range with explicit row/column
    # indices instead.
    $used = $ws.UsedRange
    $firstRow = $used.Row
    $firstCol = $used.Column
    $lastRow = $firstRow + $used.Rows.Count - 1
    $lastCol = $firstCol + $used.Columns.Count - 1

    for ($r = $firstRow; $r -le $lastRow; $r++) {
        for ($c = $firstCol; $c -le $lastCol; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $cv = $cell.Value2
            if ($null -ne $cv -and $replacements.Contains($cv)) {
                $cell.Value2 = $replacements[$cv]
            }
        }
    }

    foreach ($hl in $ws.Hyperlinks) {
        $t = $hl.TextToDisplay
        if ($null -ne $t -and $replacements.Contains($t)) {
            $hl.TextToDisplay = $replacements[$t]
        }
    }
}
